$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 545, shifting existing rows 545+ down by one.
$ws.Rows("545:545").Insert()

# Populate the newly inserted row 545 with the new record.
$ws.Range("A545").Value = 6
$ws.Range("B545").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C545").Value = "Metropolitana"
$ws.Range("D545").Value = 45218
$ws.Range("E545").Value = 13
$ws.Range("F545").Value = 100112032
$ws.Range("G545").Value = "Zapallo italiano"
$ws.Range("H545").Value = "Sin especificar"
$ws.Range("I545").Value = "Primera"
$ws.Range("J545").Value = 500
$ws.Range("K545").Value = 15000
$ws.Range("L545").Value = 17000
$ws.Range("M545").Value = 16080
$ws.Range("N545").Value = "$/caja 50 unidades"
$ws.Range("O545").Value = "Región de Arica y Parinacota"
$ws.Range("P545").Value = 322
$ws.Range("Q545").Value = 50
$ws.Range("R545").Value = "Hortaliza"
